# Auto-generated script to apply market data refresh changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 4102.6523
$ws.Cells.Item(106, 9).Value = 3597.9473
$ws.Cells.Item(106, 11).Value = 3597.9473
$ws.Cells.Item(106, 13).Value = -2966.9473
$ws.Cells.Item(129, 8).Value = 1303.1
$ws.Cells.Item(129, 10).Value = 2972.5
$ws.Cells.Item(129, 12).Value = 8917.5
$ws.Cells.Item(129, 14).Value = -18917.5
$ws.Cells.Item(137, 8).Value = 2290.2727
$ws.Cells.Item(137, 9).Value = 2286.75
$ws.Cells.Item(137, 11).Value = 6860.25
$ws.Cells.Item(137, 13).Value = -4310.25
$ws.Cells.Item(138, 8).Value = 1450
$ws.Cells.Item(138, 10).Value = 4641.4287
$ws.Cells.Item(138, 12).Value = 13924.2861
$ws.Cells.Item(138, 14).Value = -24204.2861
$ws.Cells.Item(141, 8).Value = 3929.9092
$ws.Cells.Item(141, 9).Value = 3504.1428
$ws.Cells.Item(141, 11).Value = 10512.4284
$ws.Cells.Item(141, 13).Value = -5332.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 176.125
$ws.Cells.Item(5, 9).Value = 101.8
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 101.8
$ws.Cells.Item(5, 12).Value = 300
$ws.Cells.Item(5, 13).Value = 10.2
$ws.Cells.Item(5, 14).Value = -524
$ws.Cells.Item(61, 8).Value = 5473.7856
$ws.Cells.Item(61, 9).Value = 5664.077
$ws.Cells.Item(61, 11).Value = 5664.077
$ws.Cells.Item(61, 13).Value = -5452.077
$ws.Cells.Item(74, 8).Value = 2581.2354
$ws.Cells.Item(74, 9).Value = 2502.8333
$ws.Cells.Item(74, 11).Value = 2502.8333
$ws.Cells.Item(74, 13).Value = -1628.8333
$ws.Cells.Item(77, 8).Value = 2581.2354
$ws.Cells.Item(77, 9).Value = 2502.8333
$ws.Cells.Item(77, 11).Value = 12514.1665
$ws.Cells.Item(77, 13).Value = -8146.166499999999
$ws.Cells.Item(110, 8).Value = 709.7692
$ws.Cells.Item(110, 9).Value = 692.1818
$ws.Cells.Item(110, 10).Value = 806.5
$ws.Cells.Item(110, 11).Value = 692.1818
$ws.Cells.Item(110, 12).Value = 806.5
$ws.Cells.Item(110, 13).Value = 1352.8182
$ws.Cells.Item(110, 14).Value = -4896.5
$ws.Cells.Item(122, 8).Value = 7938629.5
$ws.Cells.Item(122, 9).Value = 10102892
$ws.Cells.Item(122, 11).Value = 30308676
$ws.Cells.Item(122, 13).Value = -30306226
$ws.Cells.Item(136, 8).Value = 5473.7856
$ws.Cells.Item(136, 9).Value = 5664.077
$ws.Cells.Item(136, 11).Value = 16992.231
$ws.Cells.Item(136, 13).Value = -14442.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 176.125
$ws.Cells.Item(4, 9).Value = 101.8
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 101.8
$ws.Cells.Item(4, 12).Value = 300
$ws.Cells.Item(4, 13).Value = 13.2
$ws.Cells.Item(4, 14).Value = -530
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 4882.7144
$ws.Cells.Item(94, 9).Value = 5130.3
$ws.Cells.Item(94, 11).Value = 5130.3
$ws.Cells.Item(94, 13).Value = -4679.3
$ws.Cells.Item(99, 8).Value = 3364
$ws.Cells.Item(99, 9).Value = 2372
$ws.Cells.Item(99, 11).Value = 2372
$ws.Cells.Item(99, 13).Value = -874
$ws.Cells.Item(134, 8).Value = 3177.3125
$ws.Cells.Item(134, 9).Value = 3202.6428
$ws.Cells.Item(134, 11).Value = 9607.928400000001
$ws.Cells.Item(134, 13).Value = -7072.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6354.2
$ws.Cells.Item(31, 9).Value = 2340.8333
$ws.Cells.Item(31, 11).Value = 2340.8333
$ws.Cells.Item(31, 13).Value = -2045.8333
$ws.Cells.Item(34, 8).Value = 6354.2
$ws.Cells.Item(34, 9).Value = 2340.8333
$ws.Cells.Item(34, 11).Value = 2340.8333
$ws.Cells.Item(34, 13).Value = -2138.8333
$ws.Cells.Item(122, 8).Value = 920.2143
$ws.Cells.Item(122, 9).Value = 855.55554
$ws.Cells.Item(122, 10).Value = 1036.6
$ws.Cells.Item(122, 11).Value = 2566.66662
$ws.Cells.Item(122, 12).Value = 3109.8
$ws.Cells.Item(122, 13).Value = -116.66662
$ws.Cells.Item(122, 14).Value = -8009.799999999999
$ws.Cells.Item(132, 8).Value = 4680.3125
$ws.Cells.Item(132, 9).Value = 5084.615
$ws.Cells.Item(132, 11).Value = 15253.845
$ws.Cells.Item(132, 13).Value = -12723.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 883.24
$ws.Cells.Item(11, 9).Value = 874.2083
$ws.Cells.Item(11, 11).Value = 2622.6249
$ws.Cells.Item(11, 13).Value = -2482.6249
$ws.Cells.Item(106, 8).Value = 13464.5
$ws.Cells.Item(106, 10).Value = 13464.5
$ws.Cells.Item(106, 12).Value = 40393.5
$ws.Cells.Item(106, 14).Value = -42285.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 5372.7393
$ws.Cells.Item(80, 10).Value = 5961.846
$ws.Cells.Item(80, 12).Value = 5961.846
$ws.Cells.Item(80, 14).Value = -7957.846
$ws.Cells.Item(83, 8).Value = 5372.7393
$ws.Cells.Item(83, 10).Value = 5961.846
$ws.Cells.Item(83, 12).Value = 29809.23
$ws.Cells.Item(83, 14).Value = -39793.23
$ws.Cells.Item(113, 8).Value = 55564140
$ws.Cells.Item(113, 10).Value = 10291.667
$ws.Cells.Item(113, 12).Value = 10291.667
$ws.Cells.Item(113, 14).Value = -14631.667
$ws.Cells.Item(132, 8).Value = 2187.9736
$ws.Cells.Item(132, 9).Value = 1919.1538
$ws.Cells.Item(132, 11).Value = 5757.4614
$ws.Cells.Item(132, 13).Value = -3227.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2377.5
$ws.Cells.Item(22, 9).Value = 2881.2222
$ws.Cells.Item(22, 11).Value = 2881.2222
$ws.Cells.Item(22, 13).Value = -2586.2222
$ws.Cells.Item(27, 8).Value = 2377.5
$ws.Cells.Item(27, 9).Value = 2881.2222
$ws.Cells.Item(27, 11).Value = 2881.2222
$ws.Cells.Item(27, 13).Value = -2774.2222
$ws.Cells.Item(46, 8).Value = 3927.5715
$ws.Cells.Item(46, 10).Value = 3927.5715
$ws.Cells.Item(46, 12).Value = 3927.5715
$ws.Cells.Item(46, 14).Value = -4303.5715
$ws.Cells.Item(61, 8).Value = 7354.4443
$ws.Cells.Item(61, 9).Value = 5100
$ws.Cells.Item(61, 11).Value = 5100
$ws.Cells.Item(61, 13).Value = -4898
$ws.Cells.Item(100, 8).Value = 7916.4165
$ws.Cells.Item(100, 9).Value = 4666
$ws.Cells.Item(100, 10).Value = 8999.888999999999
$ws.Cells.Item(100, 11).Value = 4666
$ws.Cells.Item(100, 12).Value = 8999.888999999999
$ws.Cells.Item(100, 13).Value = -4125
$ws.Cells.Item(100, 14).Value = -10081.889
$ws.Cells.Item(113, 8).Value = 7354.4443
$ws.Cells.Item(113, 9).Value = 5100
$ws.Cells.Item(113, 11).Value = 5100
$ws.Cells.Item(113, 13).Value = -2930
$ws.Cells.Item(122, 8).Value = 4598.44
$ws.Cells.Item(122, 9).Value = 4647.6
$ws.Cells.Item(122, 10).Value = 4401.8
$ws.Cells.Item(122, 11).Value = 13942.8
$ws.Cells.Item(122, 12).Value = 13205.4
$ws.Cells.Item(122, 13).Value = -11492.8
$ws.Cells.Item(122, 14).Value = -18105.4
$ws.Cells.Item(132, 8).Value = 31252792
$ws.Cells.Item(132, 9).Value = 2939.238
$ws.Cells.Item(132, 10).Value = 90911600
$ws.Cells.Item(132, 11).Value = 8817.714
$ws.Cells.Item(132, 12).Value = 272734800
$ws.Cells.Item(132, 13).Value = -6287.714
$ws.Cells.Item(132, 14).Value = -272739860
$ws.Cells.Item(136, 8).Value = 21908.385
$ws.Cells.Item(136, 9).Value = 2068.75
$ws.Cells.Item(136, 11).Value = 6206.25
$ws.Cells.Item(136, 13).Value = -3656.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 14754.375
$ws.Cells.Item(62, 9).Value = 4598.4
$ws.Cells.Item(62, 10).Value = 17427
$ws.Cells.Item(62, 11).Value = 4598.4
$ws.Cells.Item(62, 12).Value = 17427
$ws.Cells.Item(62, 13).Value = -3974.4
$ws.Cells.Item(62, 14).Value = -18675
$ws.Cells.Item(63, 8).Value = 17450
$ws.Cells.Item(63, 10).Value = 17450
$ws.Cells.Item(63, 12).Value = 17450
$ws.Cells.Item(63, 14).Value = -18698
$ws.Cells.Item(65, 8).Value = 14754.375
$ws.Cells.Item(65, 9).Value = 4598.4
$ws.Cells.Item(65, 10).Value = 17427
$ws.Cells.Item(65, 11).Value = 22992
$ws.Cells.Item(65, 12).Value = 87135
$ws.Cells.Item(65, 13).Value = -19872
$ws.Cells.Item(65, 14).Value = -93375
$ws.Cells.Item(66, 8).Value = 17450
$ws.Cells.Item(66, 10).Value = 17450
$ws.Cells.Item(66, 12).Value = 52350
$ws.Cells.Item(66, 14).Value = -58590
$ws.Cells.Item(81, 8).Value = 4995.9
$ws.Cells.Item(81, 9).Value = 4425.7144
$ws.Cells.Item(81, 10).Value = 6326.3335
$ws.Cells.Item(81, 11).Value = 8851.4288
$ws.Cells.Item(81, 12).Value = 12652.667
$ws.Cells.Item(81, 13).Value = -7790.4288
$ws.Cells.Item(81, 14).Value = -14774.667
$ws.Cells.Item(84, 8).Value = 4995.9
$ws.Cells.Item(84, 9).Value = 4425.7144
$ws.Cells.Item(84, 10).Value = 6326.3335
$ws.Cells.Item(84, 11).Value = 44257.144
$ws.Cells.Item(84, 12).Value = 63263.335
$ws.Cells.Item(84, 13).Value = -38953.144
$ws.Cells.Item(84, 14).Value = -73871.33499999999
$ws.Cells.Item(100, 8).Value = 1637.5
$ws.Cells.Item(100, 9).Value = 1385
$ws.Cells.Item(100, 10).Value = 1890
$ws.Cells.Item(100, 11).Value = 2770
$ws.Cells.Item(100, 12).Value = 3780
$ws.Cells.Item(100, 13).Value = -2229
$ws.Cells.Item(100, 14).Value = -4862
$ws.Cells.Item(126, 8).Value = 1808.7368
$ws.Cells.Item(126, 9).Value = 1695.9231
$ws.Cells.Item(126, 11).Value = 5087.7693
$ws.Cells.Item(126, 13).Value = -2617.7693
